# Edit script: add Donate Now / Send E-Card / Patient Survey test rows
# (Approach / Pass Criteria / Fail Criteria columns) to the Security Testing
# worksheet, matching the new test cases added to the plan.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$donate_approach = 'We''ll test the donation system to make sure a user can make a donation online via the donation form and then be emailed a receipt. The user will complete the form using the name Dave Davidson, with dryden_dave_test@gmail.com as the email address and a test credit card account. The user will submit the donation form and test for success. The user will then check the email account to ensure the receipt arrives. The user will also ensure the data has been inserted into the database.'
$donate_pass = 'Proper validation of the form, the payment is secure and the transaction only occurs if all validation methods are met, the receipt is emailed to the donor, and all of the information is inserted into the database. The feature is responsive and works on all browsers.'
$donate_fail = 'The transaction doesn''t go through, the transaction goes through multiple times, validation doesn''t occur, transaction not secure, information doesn''t get stored in the database, no receipt is emailed. The page isn''t responsive and/or doesn''t work on all browsers.'
$ecard_approach = 'We''ll test the e-card template that will allow users to send an e-card to a patient at the hospital. We''ll test by having a user enter a name and an email adress of dryden_dave_test@gmail.com as the patients'' email. The user will then fill out a greeting and customize the card with their choice of background and font. The user will then send the card and be greeted with a success message. The user will then check the email to see if the card was successfully sent and appears as designed. The user will then test by using improper information and blank fields to ensure all error messages appear and ensure it does not send.'
$ecard_pass = 'The user is successfully able to send an e-card they designed, the user receives a success message upon sending the card, proper validation occurs and errors occur when fields aren''t filled in properly or left blank. The feature is responsive and works on all browsers.'
$ecard_fail = 'Validation doesn''t work as designed, e-card isn''t sent, e-card isn''t as designed, page is not responsive and/or doesn''t work on all browsers.'
$survey_approach = 'We''ll test the  control of the patient survey to make sure the steps allow the user to go forward and backwards appropriately, proper validation controls and messages occur, and that the results are stored in the database. The user will go through the survey and then submit, and test for a success message. After a successful message, we will verify the answers were submitted into the database. The user will do the survey again, but this time leave questions unanswered and submit. Testing if the appropriate error messages appear. '
$survey_pass = 'The wizard control functions as designed with the buttons and controls, as well as displaying the proper error messages and the results are successfully inserted into the database. The user receives a thank you message thanking them for their feedback upon successful completion of the survey.The feature is responsive and works on all browsers.'
$survey_fail = 'The wizard control does not provide proper functionality, accepts invalid input, does not display error messages, and does not store the data in the database, page is not responsive and/or doesn''t work on all browsers.'

# Row 15 - Donate Now
$ws.Range("G15").Value = $donate_approach
$ws.Range("H15").Value = $donate_pass
$ws.Range("I15").Value = $donate_fail
$ws.Range("G15:I15").WrapText = $true
$ws.Rows.Item(15).RowHeight = 195

# Row 16 - Send E-Card
$ws.Range("G16").Value = $ecard_approach
$ws.Range("H16").Value = $ecard_pass
$ws.Range("I16").Value = $ecard_fail
$ws.Range("G16:I16").WrapText = $true
$ws.Rows.Item(16).RowHeight = 240

# Row 17 - Patient Survey
$ws.Range("G17").Value = $survey_approach
$ws.Range("H17").Value = $survey_pass
$ws.Range("I17").Value = $survey_fail
$ws.Range("G17:I17").WrapText = $true
$ws.Rows.Item(17).RowHeight = 195

# Update the saved view/selection state to reflect scrolling down to the
# newly added rows and selecting the new approach/pass/fail block.
[void]$ws.Activate()
[void]$ws.Range("G15:I17").Select()
